# Replace every column-F "(0, 0)" placeholder string with "(nan, nan)" on
# every worksheet in the workbook. This mirrors an upstream documentation
# fix where the coefficient-extraction helper started reporting missing
# "EXP Start Point" values as NaN pairs instead of zero pairs.

$wb = $excel.ActiveWorkbook

$changed = 0

foreach ($ws in $wb.Worksheets) {
    $usedRange = $ws.UsedRange
    $lastRow = $usedRange.Rows.Count

    for ($r = 1; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 6)
        $val = $cell.Value2
        if ($val -eq "(0, 0)") {
            $cell.Value = "(nan, nan)"
            $changed++
        }
    }
}

Write-Output "cells changed: $changed"
